$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The averaged-intensity scan was re-run to include the new spiral sampling
# schemes. "Gaussian-Quadrature" (previously the last method, row 16) and
# three new "Spiral-*" schemes are inserted right after
# "Ring Perpendicular to TD", pushing the remaining rotation/hex-grid
# methods down by three rows. The sheet grows from 14 methods (rows 3-16)
# to 17 methods (rows 3-19).

$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("B16").Value = "Rotation-60detTilt"

# New rows 17-19 hold the methods that used to live in rows 13-15
# (HexGrid-90degTilt5degRes, HexGrid-90degTilt22p5degRes,
# HexGrid-60degTilt5degRes), each with the same "1" intensity values
# across every HKL column (C:P) as every other row.
$newRows = @(
    @{ Row = 17; Index = 15; Name = "HexGrid-90degTilt5degRes" },
    @{ Row = 18; Index = 16; Name = "HexGrid-90degTilt22p5degRes" },
    @{ Row = 19; Index = 17; Name = "HexGrid-60degTilt5degRes" }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value = $nr.Index
    $ws.Cells.Item($r, 2).Value = $nr.Name
    for ($col = 3; $col -le 16; $col++) {
        $ws.Cells.Item($r, $col).Value = 1
    }
    # Column A carries the same bold/centered/bordered "index" style as
    # every other row (e.g. row 16) - copy it across instead of hard
    # coding font/border constants.
    $ws.Range("A16").Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false
